$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "instructor"

# Data rows 8-19 (columns A,B,C,D,F,G,H,I; E intentionally left blank)
$rows = @(
  @{ r=8;  A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start A small MMA studio in my garage. I already have a small clientele base who is interested and instructors who would work for me. "; F="1000-1500"; G="N/A"; H="MMA"; I="No" },
  @{ r=9;  A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling"; I="Yes" },
  @{ r=10; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=11; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=12; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=13; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=14; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=15; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=16; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="I want to start a big wrestling gym in a commercial area."; F="1000-1500"; G="N/A"; H="Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=17; A="Daniel"; B="dansem@gnail.com"; C="State college"; D="Gym"; F="500-1000"; G=""; H="Muay Thai"; I="Yes" },
  @{ r=18; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="MMA"; F="1000-1500"; G=""; H="Muay Thai, Wrestling, Jiu-Jitsu"; I="Yes" },
  @{ r=19; A="Frank Roman Bevivino"; B="frankbevivino@gmail.com"; C="State College"; D="MMA"; F="1000-1500"; G=""; H="Muay Thai, Wrestling, Jiu-Jitsu"; I="Yes" }
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  # column E (disciplines[]) is intentionally left blank for these rows
  $ws.Cells.Item($r, 6).Value = $row.F
  if ($row.G -ne "") {
    $ws.Cells.Item($r, 7).Value = $row.G
  }
  $ws.Cells.Item($r, 8).Value = $row.H
  $ws.Cells.Item($r, 9).Value = $row.I
}
